$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.714.51"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.232.76"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.10%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.48"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.21%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.232.68"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.73"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.504"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.46%  "

$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.97"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.764.18"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.757.31"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.234.94"
$ws.Range("D18").ClearFormats()

$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.52"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.62"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.04"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.16"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("E29").Value = "  +2.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +36.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.94"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.22"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("E35").Value = "  -4.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.50"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.65"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "502.94"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0775"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +14.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.11"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.91%  "

$ws.Range("E41").Value = "  -0.71%  "

$ws.Range("E42").Value = "  +6.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.73"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.80%  "

$ws.Range("E44").Value = "  -1.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.913.98"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.19"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.87%  "

$ws.Range("E48").Value = "  +4.15%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("E50").Value = "  -0.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.14%  "
